$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: new "Python" entry -------------------------------------------------
# Shared strings must be appended in this exact order so they land at
# indices 24, 25, 26 (Python / path / help text).
$ws.Range("A9").Value = "Python"
$ws.Range("B9").Value = "C:/PROGRA~1/Python310/python.exe"
$ws.Range("C9").Value = "Location of the Python executable to use - some parts of these workflows will run some python code."

# --- Styling for B9 / C9 --------------------------------------------------------
# B9 should look like the existing "path" column style (italic Consolas) and
# C9 like the existing "help" column style (italic Calibri). We seed each
# cell from a base that already has italic switched on (so only one font
# property actually needs to change), using mismatched source/target fonts
# so that a genuinely new font/style entry is created rather than reusing an
# existing one.
$ws.Range("C2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("B9").Font.Name = "Consolas"
$ws.Range("C9").Font.Name = "Calibri"

# --- Selection -------------------------------------------------------------
$ws.Range("B12").Select()
